$d = $word.ActiveDocument

# 1. Title replacement (appears twice: Heading1 title and bold run near end)
$d.Content.Find.Execute("Play Book of Gold: Symbol Choice for Free - Review", $true, $false, $false, $false, $false, $true, 1, $false, "Play Book of Gold: Symbol Choice Free | Slot Game Review", 2)

# 2. "What we like" bullet updates
$d.Content.Find.Execute("Simple and intuitive gameplay", $true, $false, $false, $false, $false, $true, 1, $false, "Simple and intuitive gameplay mechanics", 2)

$d.Content.Find.Execute("Chance to win significant amounts through Free Spins", $true, $false, $false, $false, $false, $true, 1, $false, "Chance to win significant amounts through Free Spins mode", 2)

$d.Content.Find.Execute("Well-executed Ancient Egyptian theme", $true, $false, $false, $false, $false, $true, 1, $false, "Decent theoretical return to player (RTP)", 2)

$d.Content.Find.Execute("Golden book serves as both wild and scatter symbol", $true, $false, $false, $false, $false, $true, 1, $false, "Well-done Egyptian-themed graphics and symbols", 2)

# 3. "What we don't like" bullet updates
$d.Content.Find.Execute("Average RTP", $true, $false, $false, $false, $false, $true, 1, $false, "Lack of groundbreaking features", 2)

$d.Content.Find.Execute("Graphics and features not particularly impressive", $true, $false, $false, $false, $false, $true, 1, $false, "Graphics could be more impressive", 2)

# 4. Final italic summary paragraph
$d.Content.Find.Execute("Read our review of Book of Gold: Symbol Choice. Play for free and enjoy the traditional gameplay and Egyptian theme. Try your luck with Free Spins!", $true, $false, $false, $false, $false, $true, 1, $false, "Read our review of Book of Gold: Symbol Choice, a slot game with a chance to win big. Play for free today!", 2)
